$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 39589556
$ws.Range("C1").Value = 39524556
$ws.Range("D1").Value = "Было Katta Doimiy 100, Стало: Katta Doimiy 40"
$ws.Range("E1").Value = "2024-10-24 23:33:05"
